$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "310.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.30%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.32"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.48%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.126"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.22%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07697"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.23%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.626"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.57%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9234"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.88%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.468"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.57%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1212"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "19.35%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1832"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.69%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09149"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.54%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04215"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-1.35%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1050"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.52%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001253"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.58%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005708"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.31%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.351"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.03%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.312"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.97%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.92%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.936"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "5.99%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.30%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.90%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04048"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-3.07%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001262"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.81%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004103"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.77%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001269"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.33%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "24.54%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02470"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.29%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05267"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.55%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007836"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.81%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1316"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.58%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006789"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.81%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001843"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.71%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008187"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.90%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3096"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-6.94%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006721"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "5.62%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.30%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "2,790.18%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004094"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-7.10%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002097"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.30%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001997"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.30%"
